# TORCH narrative + figures
# 1) Update the fixed "datetimeFigureOut" field text (24/06/2020 -> 27/06/2020)
#    on the Slide Master and every Slide Layout's date placeholder.
# 2) Update the slide title text (remove curly quotes around Torch, upper-case it).
# 3) Resize/reposition the title textbox slightly (per author's figure tweak).

$p = $ppt.ActivePresentation

$oldDate = "24/06/2020"
$newDate = "27/06/2020"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# Slide 1 title text + textbox geometry.
$s = $p.Slides.Item(1)
$title = $s.Shapes.Item(11)
$title.TextFrame.TextRange.Text = "Example E7 – Part 2: TORCH super panel  including tests and sub-panels"
$title.Left = 10.9772
$title.Width = 949.0227559055118
